$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize cell formatting first (remove the one-off black-font style applied to
#     C2, D2, E2, C3, D3, E3) by copying the plain bordered format (style used by B2)
#     onto those cells. This matches the original/older template formatting where
#     every data cell in the table body shares the same simple bordered style. ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2:E2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3:E3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Header row (row 1): restore the original Lob / Process / Product Name column order ---
$ws.Range("F1").Value = "Product Name"
$ws.Range("G1").Value = "Lob"
$ws.Range("H1").Value = "Process"

# --- Row 2 data ---
$ws.Range("A2").Value = 45436
$ws.Range("B2").Value = 1213286
$ws.Range("C2").Value = "SIPL0005"
$ws.Range("D2").Value = "SIPL0004"
$ws.Range("E2").Value = "WFG Title"
$ws.Range("F2").Value = "Full Search"
$ws.Range("G2").Value = "Title"
$ws.Range("H2").Value = "Search"
$ws.Range("I2").Value = "FL"
$ws.Range("J2").Value = "Clay"
$ws.Range("K2").Value = "FLClay"
$ws.Range("L2").Value = "WIP"
$ws.Range("M2").Value = "Search(T1) "

# --- Row 3 data ---
$ws.Range("A3").Value = 45439
$ws.Range("B3").Value = 2193289
$ws.Range("C3").Value = "SIPL0005"
$ws.Range("D3").Value = "SIPL0004"
$ws.Range("E3").Value = "WFG Title"
$ws.Range("F3").Value = "Current Owner Search"
$ws.Range("G3").Value = "Title"
$ws.Range("H3").Value = "Search"
$ws.Range("I3").Value = "FL"
$ws.Range("J3").Value = "Clay"
$ws.Range("K3").Value = "FLClay"
$ws.Range("L3").Value = "WIP"
$ws.Range("M3").Value = "Search(T2)"

# --- Column widths: widen Emp ID-Order Assigned (C) and make F:H uniform width,
#     matching F/G's existing width so H loses its old bestFit/auto width. ---
$ws.Columns.Item(3).ColumnWidth = 36.2
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# --- Selection cursor position ---
$ws.Range("I13").Select() | Out-Null
